$wb = $excel.ActiveWorkbook

function Set-HeaderBorderStyle($cell, [bool]$includeRight) {
    # Start from a clean slate so the new style does not inherit the bold
    # font / full border / alignment of the merged header cell (s=1).
    $cell.ClearFormats()
    $top = $cell.Borders.Item(8)     # xlEdgeTop
    $top.LineStyle = 1               # xlContinuous
    $top.Weight = 2                  # xlThin
    if ($includeRight) {
        $right = $cell.Borders.Item(10)  # xlEdgeRight
        $right.LineStyle = 1
        $right.Weight = 2
    }
    $bottom = $cell.Borders.Item(9)  # xlEdgeBottom
    $bottom.LineStyle = 1
    $bottom.Weight = 2
}

# ---------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-HeaderBorderStyle $ws1.Range("C1") $false
Set-HeaderBorderStyle $ws1.Range("D1") $true

$ws1.Range("C2").Value = "approach"

$ws1.Range("D4").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("D12").Value = 0

# ---------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-HeaderBorderStyle $ws2.Range("C1") $false
Set-HeaderBorderStyle $ws2.Range("D1") $true
Set-HeaderBorderStyle $ws2.Range("F1") $false
Set-HeaderBorderStyle $ws2.Range("G1") $true

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

$ws2.Range("G5").ClearContents()
